# Month-append update: shift last month's "E" column figures into "B" (the
# prior-month reference column) and populate "E" with the new month's
# figures. Row/column labels (A, D) are untouched; only data columns move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: month headers -------------------------------------------------
# B1 used to hold the literal text "Apr"; it now becomes a real date serial
# (matching E1's format) for the month that just rolled off, and E1 advances
# to the new month.
$ws.Range("B1").NumberFormat = "mmm-yy"
$ws.Range("B1").Value = 45413
$ws.Range("E1").Value = 45413

# --- Left table (A/B): previous-month column now carries what used to be
# this month's figures -----------------------------------------------------
$ws.Range("B3").Value = 3429675
$ws.Range("B4").Value = 7010
$ws.Range("B6").Value = 4168133
$ws.Range("B8").Value = 51596
$ws.Range("B10").Value = 100781.58
$ws.Range("B11").Value = 23836
$ws.Range("B14").ClearContents()
$ws.Range("B15").Value = 532255
$ws.Range("B16").ClearContents()
$ws.Range("B18").Value = 22550
$ws.Range("B19").Value = 204625
$ws.Range("B21").Value = 9136834.5800000001
$ws.Range("B23").Value = 4227041
$ws.Range("B25").Value = 616661
$ws.Range("B26").Value = 4168133
$ws.Range("B27").Value = 9011835
$ws.Range("B29").Value = 9136835
$ws.Range("B30").Value = 4784794
$ws.Range("B32").Value = 8271.5
$ws.Range("B33").Value = 52769.8
$ws.Range("B34").Value = 63524
$ws.Range("B35").Value = 1229816
$ws.Range("B36").Value = 10366651
$ws.Range("B38").Value = 9136834.5800000001

# --- Right table (D/E): refreshed with the new month's figures ------------
$ws.Range("E2").Value = 259749
$ws.Range("E3").Value = 629629
$ws.Range("E4").Value = 4168133
$ws.Range("E5").Value = 744437
$ws.Range("E6").Value = 204625
$ws.Range("E7").Value = 100781.58
$ws.Range("E8").Value = 532255
$ws.Range("E9").Value = 22550
$ws.Range("E10").Value = 51596
$ws.Range("E13").Value = 30846
$ws.Range("E14").Value = 7215974.5800000001
$ws.Range("E15").Value = 3610.26
$ws.Range("E16").Value = 172490.69
$ws.Range("E17").Value = 399870
$ws.Range("E18").Value = 575970.94999999995
$ws.Range("E19").Value = 7791945.5300000003
$ws.Range("E22").Value = 863710
$ws.Range("E23").Value = 407150
$ws.Range("E26").Value = 368463.4
$ws.Range("E27").Value = 50771.509999999995
$ws.Range("E28").Value = 10384
$ws.Range("E29").Value = 224226
$ws.Range("E30").Value = 2574704.9099999997
$ws.Range("E31").Value = 10366650.439999999
$ws.Range("E32").Value = 10366651
$ws.Range("E33").Value = -0.56000000052154064
$ws.Range("E34").Value = 2278088
$ws.Range("E35").Value = 3429675
$ws.Range("E37").Value = 9136834.5800000001
$ws.Range("E38").Value = 1229815.8599999999
$ws.Range("E39").Value = 10366650.439999999
$ws.Range("E40").Value = 9136834.5800000001
$ws.Range("E41").Value = 1229815.8599999999
$ws.Range("E42").Value = 10366650.439999999
$ws.Range("E44").Value = 2641134.5300000003
$ws.Range("E45").Value = 5057511
$ws.Range("E46").Value = 93300
$ws.Range("E47").Value = 2449704.9099999997
$ws.Range("E49").Value = 10366650.439999999

# --- View: scroll back to A1 (drop the stale frozen-ish topLeftCell) and
# select the newly-updated May figures ------------------------------------
$ws.Range("E22:E49").Select()
